# Todo.xlsx update: add three new TODO rows to Sheet1 and reposition the
# active window/selection, matching the author's latest commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 29-31 in column C -------------------------------------------------
# Row 29 and 30 get the same "white / Background 1" fill highlight that the
# last completed-looking items further up the sheet use; row 31 is left with
# the sheet's default (no) style, matching a freshly typed, not-yet-styled row.

$ws.Range("C29").Value = "Line numbers"
$ws.Range("C29").Interior.ThemeColor = 2   # xlThemeColorLight1 -> theme="0" (Background 1)
$ws.Range("C29").Interior.TintAndShade = 0

$ws.Range("C30").Value = "Fix usage in loops…"
$ws.Range("C30").Interior.ThemeColor = 2   # xlThemeColorLight1 -> theme="0" (Background 1)
$ws.Range("C30").Interior.TintAndShade = 0

$ws.Range("C31").Value = "Remove dead declarations…"

# --- View state -------------------------------------------------------------
# Move the active selection to the next empty row under the new entries, as
# happens naturally once the user finishes typing the last item.
$ws.Range("C32").Select()

# Reposition/resize the workbook window to match the saved view (best effort;
# harmless if the host does not persist window chrome).
$win = $excel.ActiveWindow
$win.Left = -60110
$win.Top = -19380
$win.Width = 21820
$win.Height = 38020
